# Appends 10 new match rows (rows 12-21) to the "Priyam Garg" batting log sheet.
# Source: webScrapping/espn_scrapper/IPL/Sunrisers Hyderabad/Priyam Garg .xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("A12").Value = " Dubai (DSC)"
$ws.Range("B12").Value = " October 13 2020"
$ws.Range("C12").Value = "Super Kings won by 20 runs"
$ws.Range("D12").Value = "Sunrisers Hyderabad"
$ws.Range("E12").Value = "Chennai Super Kings"
$ws.Range("F12").Value = "Priyam Garg "
$ws.Range("G12").Value = "'16"
$ws.Range("H12").Value = "'18"
$ws.Range("I12").Value = "'1"
$ws.Range("J12").Value = "'0"
$ws.Range("K12").Value = "'88.88"

# Row 13
$ws.Range("A13").Value = " Abu Dhabi"
$ws.Range("B13").Value = " October 18 2020"
$ws.Range("C13").Value = "Match tied (KKR won the one-over eliminator)"
$ws.Range("D13").Value = "Sunrisers Hyderabad"
$ws.Range("E13").Value = "Kolkata Knight Riders"
$ws.Range("F13").Value = "Priyam Garg "
$ws.Range("G13").Value = "'4"
$ws.Range("H13").Value = "'7"
$ws.Range("I13").Value = "'0"
$ws.Range("J13").Value = "'0"
$ws.Range("K13").Value = "'57.14"

# Row 14
$ws.Range("A14").Value = " Dubai (DSC)"
$ws.Range("B14").Value = " October 02 2020"
$ws.Range("C14").Value = "Sunrisers won by 7 runs"
$ws.Range("D14").Value = "Sunrisers Hyderabad"
$ws.Range("E14").Value = "Chennai Super Kings"
$ws.Range("F14").Value = "Priyam Garg "
$ws.Range("G14").Value = "'51"
$ws.Range("H14").Value = "'26"
$ws.Range("I14").Value = "'6"
$ws.Range("J14").Value = "'1"
$ws.Range("K14").Value = "'196.15"

# Row 15
$ws.Range("A15").Value = " Sharjah"
$ws.Range("B15").Value = " October 04 2020"
$ws.Range("C15").Value = "Mumbai won by 34 runs"
$ws.Range("D15").Value = "Sunrisers Hyderabad"
$ws.Range("E15").Value = "Mumbai Indians"
$ws.Range("F15").Value = "Priyam Garg "
$ws.Range("G15").Value = "'8"
$ws.Range("H15").Value = "'7"
$ws.Range("I15").Value = "'0"
$ws.Range("J15").Value = "'0"
$ws.Range("K15").Value = "'114.28"

# Row 16
$ws.Range("A16").Value = " Abu Dhabi"
$ws.Range("B16").Value = " November 06 2020"
$ws.Range("C16").Value = "Sunrisers won by 6 wickets (with 2 balls remaining)"
$ws.Range("D16").Value = "Sunrisers Hyderabad"
$ws.Range("E16").Value = "Royal Challengers Bangalore"
$ws.Range("F16").Value = "Priyam Garg "
$ws.Range("G16").Value = "'7"
$ws.Range("H16").Value = "'14"
$ws.Range("I16").Value = "'0"
$ws.Range("J16").Value = "'0"
$ws.Range("K16").Value = "'50.00"

# Row 17
$ws.Range("A17").Value = " Dubai (DSC)"
$ws.Range("B17").Value = " September 21 2020"
$ws.Range("C17").Value = "RCB won by 10 runs"
$ws.Range("D17").Value = "Sunrisers Hyderabad"
$ws.Range("E17").Value = "Royal Challengers Bangalore"
$ws.Range("F17").Value = "Priyam Garg "
$ws.Range("G17").Value = "'12"
$ws.Range("H17").Value = "'13"
$ws.Range("I17").Value = "'1"
$ws.Range("J17").Value = "'0"
$ws.Range("K17").Value = "'92.30"

# Row 18
$ws.Range("A18").Value = " Abu Dhabi"
$ws.Range("B18").Value = " November 08 2020"
$ws.Range("C18").Value = "Capitals won by 17 runs"
$ws.Range("D18").Value = "Sunrisers Hyderabad"
$ws.Range("E18").Value = "Delhi Capitals"
$ws.Range("F18").Value = "Priyam Garg "
$ws.Range("G18").Value = "'17"
$ws.Range("H18").Value = "'12"
$ws.Range("I18").Value = "'0"
$ws.Range("J18").Value = "'2"
$ws.Range("K18").Value = "'141.66"

# Row 19
$ws.Range("A19").Value = " Dubai (DSC)"
$ws.Range("B19").Value = " October 11 2020"
$ws.Range("C19").Value = "Royals won by 5 wickets (with 1 ball remaining)"
$ws.Range("D19").Value = "Sunrisers Hyderabad"
$ws.Range("E19").Value = "Rajasthan Royals"
$ws.Range("F19").Value = "Priyam Garg "
$ws.Range("G19").Value = "'15"
$ws.Range("H19").Value = "'8"
$ws.Range("I19").Value = "'1"
$ws.Range("J19").Value = "'1"
$ws.Range("K19").Value = "'187.50"

# Row 20
$ws.Range("A20").Value = " Dubai (DSC)"
$ws.Range("B20").Value = " October 24 2020"
$ws.Range("C20").Value = "Kings XI won by 12 runs"
$ws.Range("D20").Value = "Sunrisers Hyderabad"
$ws.Range("E20").Value = "Kings XI Punjab"
$ws.Range("F20").Value = "Priyam Garg "
$ws.Range("G20").Value = "'3"
$ws.Range("H20").Value = "'5"
$ws.Range("I20").Value = "'0"
$ws.Range("J20").Value = "'0"
$ws.Range("K20").Value = "'60.00"

# Row 21
$ws.Range("A21").Value = " Dubai (DSC)"
$ws.Range("B21").Value = " October 08 2020"
$ws.Range("C21").Value = "Sunrisers won by 69 runs"
$ws.Range("D21").Value = "Sunrisers Hyderabad"
$ws.Range("E21").Value = "Kings XI Punjab"
$ws.Range("F21").Value = "Priyam Garg "
$ws.Range("G21").Value = "'0"
$ws.Range("H21").Value = "'1"
$ws.Range("I21").Value = "'0"
$ws.Range("J21").Value = "'0"
$ws.Range("K21").Value = "'0.00"

